# Auto-generated edit script applying numeric corrections to the Leviathan_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ----- ALC sheet -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1076.7428
$ws.Range("I15").Value = 1076.7428
$ws.Range("K15").Value = 3230.2284
$ws.Range("M15").Value = -3061.2284
$ws.Range("H17").Value = 477891.53
$ws.Range("J17").Value = 501706.16
$ws.Range("L17").Value = 1505118.48
$ws.Range("N17").Value = -1505454.48
$ws.Range("H18").Value = 4580
$ws.Range("I18").Value = 4580
$ws.Range("K18").Value = 4580
$ws.Range("M18").Value = -4296
$ws.Range("H32").Value = 3912.0557
$ws.Range("I32").Value = 4011.9092
$ws.Range("J32").Value = 3755.1428
$ws.Range("K32").Value = 4011.9092
$ws.Range("L32").Value = 3755.1428
$ws.Range("M32").Value = -3685.9092
$ws.Range("N32").Value = -4407.1428
$ws.Range("H38").Value = 352.29413
$ws.Range("I38").Value = 352.29413
$ws.Range("K38").Value = 1056.88239
$ws.Range("M38").Value = -684.88239
$ws.Range("H97").Value = 2480.5
$ws.Range("J97").Value = 2480.5
$ws.Range("L97").Value = 7441.5
$ws.Range("N97").Value = -8433.5
$ws.Range("H132").Value = 2090.255
$ws.Range("I132").Value = 832.0625
$ws.Range("K132").Value = 2496.1875
$ws.Range("M132").Value = 33.8125
$ws.Range("H135").Value = 1570.5217
$ws.Range("I135").Value = 1451.0625
$ws.Range("J135").Value = 1843.5714
$ws.Range("K135").Value = 13059.5625
$ws.Range("L135").Value = 16592.1426
$ws.Range("M135").Value = -10524.5625
$ws.Range("N135").Value = -21662.1426
$ws.Range("H137").Value = 3635.3635
$ws.Range("I137").Value = 2698.4285
$ws.Range("K137").Value = 8095.2855
$ws.Range("M137").Value = -5545.2855
$ws.Range("H138").Value = 3435.611
$ws.Range("J138").Value = 3863.6956
$ws.Range("L138").Value = 11591.0868
$ws.Range("N138").Value = -21871.0868

# ----- ARM sheet -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7405.2
$ws.Range("J46").Value = 7472
$ws.Range("L46").Value = 7472
$ws.Range("N46").Value = -8110
$ws.Range("H74").Value = 1624.6072
$ws.Range("I74").Value = 938.15
$ws.Range("J74").Value = 3340.75
$ws.Range("K74").Value = 938.15
$ws.Range("L74").Value = 3340.75
$ws.Range("M74").Value = -64.14999999999998
$ws.Range("N74").Value = -5088.75
$ws.Range("H77").Value = 1624.6072
$ws.Range("I77").Value = 938.15
$ws.Range("J77").Value = 3340.75
$ws.Range("K77").Value = 4690.75
$ws.Range("L77").Value = 16703.75
$ws.Range("M77").Value = -322.75
$ws.Range("N77").Value = -25439.75
$ws.Range("H132").Value = 3347.3635
$ws.Range("I132").Value = 3037.3333
$ws.Range("K132").Value = 9111.999899999999
$ws.Range("M132").Value = -6581.999899999999
$ws.Range("H138").Value = 119997.2
$ws.Range("J138").Value = 119997.2
$ws.Range("L138").Value = 119997.2
$ws.Range("N138").Value = -130277.2

# ----- BSM sheet -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H20").Value = 25007764
$ws.Range("I20").Value = 33343170
$ws.Range("K20").Value = 33343170
$ws.Range("M20").Value = -33342923
$ws.Range("H134").Value = 108301.29
$ws.Range("I134").Value = 121142.28
$ws.Range("J134").Value = 1293
$ws.Range("K134").Value = 363426.84
$ws.Range("L134").Value = 3879
$ws.Range("M134").Value = -360891.84
$ws.Range("N134").Value = -8949

# ----- CRP sheet -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1236.8889
$ws.Range("J16").Value = 1385.3334
$ws.Range("L16").Value = 1385.3334
$ws.Range("N16").Value = -1959.3334
$ws.Range("H31").Value = 3833.2888
$ws.Range("I31").Value = 2474.4644
$ws.Range("K31").Value = 2474.4644
$ws.Range("M31").Value = -2179.4644
$ws.Range("H34").Value = 3833.2888
$ws.Range("I34").Value = 2474.4644
$ws.Range("K34").Value = 2474.4644
$ws.Range("M34").Value = -2272.4644
$ws.Range("H51").Value = 1000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 20746.25
$ws.Range("J68").Value = 20746.25
$ws.Range("L68").Value = 20746.25
$ws.Range("N68").Value = -22244.25
$ws.Range("H71").Value = 20746.25
$ws.Range("J71").Value = 20746.25
$ws.Range("L71").Value = 62238.75
$ws.Range("N71").Value = -69726.75
$ws.Range("H74").Value = 58216.668
$ws.Range("J74").Value = 58216.668
$ws.Range("L74").Value = 58216.668
$ws.Range("N74").Value = -59964.668
$ws.Range("H77").Value = 58216.668
$ws.Range("J77").Value = 58216.668
$ws.Range("L77").Value = 174650.004
$ws.Range("N77").Value = -183386.004
$ws.Range("H105").Value = 1496.4412
$ws.Range("I105").Value = 1354.4482
$ws.Range("J105").Value = 2320
$ws.Range("K105").Value = 1354.4482
$ws.Range("L105").Value = 2320
$ws.Range("M105").Value = 392.5518
$ws.Range("N105").Value = -5814
$ws.Range("H107").Value = 2193.6667
$ws.Range("I107").Value = 1687
$ws.Range("J107").Value = 2869.2222
$ws.Range("K107").Value = 1687
$ws.Range("L107").Value = 2869.2222
$ws.Range("M107").Value = 233
$ws.Range("N107").Value = -6709.2222
$ws.Range("H113").Value = 1236.8889
$ws.Range("J113").Value = 1385.3334
$ws.Range("L113").Value = 1385.3334
$ws.Range("N113").Value = -5725.3334
$ws.Range("H132").Value = 3176
$ws.Range("I132").Value = 3066.2
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 9198.599999999999
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -6668.599999999999
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 2339.7932
$ws.Range("I134").Value = 1993.0754
$ws.Range("K134").Value = 5979.2262
$ws.Range("M134").Value = -3444.2262
$ws.Range("H141").Value = 237857.2
$ws.Range("J141").Value = 237857.2
$ws.Range("L141").Value = 237857.2
$ws.Range("N141").Value = -248217.2

# ----- CUL sheet -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2778519.5
$ws.Range("I11").Value = 3750250.5
$ws.Range("K11").Value = 11250751.5
$ws.Range("M11").Value = -11250611.5
$ws.Range("H55").Value = 9696022
$ws.Range("J55").Value = 25005298
$ws.Range("L55").Value = 75015894
$ws.Range("N55").Value = -75016248
$ws.Range("H68").Value = 1772.3334
$ws.Range("I68").Value = 1284
$ws.Range("K68").Value = 3852
$ws.Range("M68").Value = -3041
$ws.Range("H71").Value = 1772.3334
$ws.Range("I71").Value = 1284
$ws.Range("K71").Value = 11556
$ws.Range("M71").Value = -7500
$ws.Range("H132").Value = 2730.2
$ws.Range("J132").Value = 3233
$ws.Range("L132").Value = 29097
$ws.Range("N132").Value = -34157

# ----- GSM sheet -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4997.6665
$ws.Range("I132").Value = 3098.9092
$ws.Range("J132").Value = 8795.182000000001
$ws.Range("K132").Value = 9296.7276
$ws.Range("L132").Value = 26385.546
$ws.Range("M132").Value = -6766.7276
$ws.Range("N132").Value = -31445.546
$ws.Range("H134").Value = 59034
$ws.Range("J134").Value = 59034
$ws.Range("L134").Value = 177102
$ws.Range("N134").Value = -182172

# ----- LTW sheet -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4378.364
$ws.Range("I22").Value = 4444
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 4444
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -4149
$ws.Range("N22").Value = -3590
$ws.Range("H27").Value = 4378.364
$ws.Range("I27").Value = 4444
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 4444
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -4337
$ws.Range("N27").Value = -3214
$ws.Range("H62").Value = 23124.5
$ws.Range("J62").Value = 23124.5
$ws.Range("L62").Value = 23124.5
$ws.Range("N62").Value = -24372.5
$ws.Range("H65").Value = 23124.5
$ws.Range("J65").Value = 23124.5
$ws.Range("L65").Value = 69373.5
$ws.Range("N65").Value = -75613.5
$ws.Range("H132").Value = 4356.4707
$ws.Range("I132").Value = 3932.6206
$ws.Range("K132").Value = 11797.8618
$ws.Range("M132").Value = -9267.861800000001

# ----- WVR sheet -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3755.5557
$ws.Range("I132").Value = 2995.348
$ws.Range("K132").Value = 8986.044
$ws.Range("M132").Value = -6456.044

